$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 96.95999200548977
$ws.Range("D4").Value = 51.40427538286895
$ws.Range("E4").Value = 0.5301596495589489
$ws.Range("F4").Value = 1.886224273823784
$ws.Range("G4").Value = 141.1196250437642
$ws.Range("H4").Value = 0.02433365207616589
$ws.Range("I4").Value = 0.3097206776392341
$ws.Range("J4").Value = 4.825945775955915
$ws.Range("K4").Value = 1.434611532604322
$ws.Range("L4").Value = 5.622603700961918
$ws.Range("M4").Value = 0.2179608714068308
$ws.Range("N4").Value = 0.1999594387507386
$ws.Range("O4").Value = 0.0008709411613381235
$ws.Range("P4").Value = 28.30095673212782
$ws.Range("Q4").Value = 5.496810862794518
$ws.Range("R4").Value = 51.23543111188337
$ws.Range("S4").Value = -1410.978289566236
$ws.Range("T4").Value = -0.04337708201092028
$ws.Range("U4").Value = -3.096335835231002
$ws.Range("V4").Value = -19.95850102743134
$ws.Range("W4").Value = -50.72922614682466
$ws.Range("X4").Value = 36.88931578584015
